$d = $word.ActiveDocument

$d.Content.Find.Execute("ZA-Information / Zentralarchiv für Empirische Sozialforschung", $true, $false, $false, $false, $false, $true, 1, $false, "ZUMA Nachrichten", 2)

$d.Content.Find.Execute("1997", $true, $false, $false, $false, $false, $true, 1, $false, "2009", 2)

$d.Content.Find.Execute("An Evaluation of Object-Oriented DBMS Developments: 1994 Edition.", $true, $false, $false, $false, $false, $true, 1, $false, "Integrating Object-Oriented Applications and Middleware with Relational Databases.", 2)

$d.Content.Find.Execute("Hochschule Aalen", $true, $false, $false, $false, $false, $true, 1, $false, "Ostbayerische Technische Hochschule Amberg-Weiden", 2)

$d.Content.Find.Execute("Buelowstrasse 80", $true, $false, $false, $false, $false, $true, 1, $false, "Albrechtstrasse 84", 2)
